# Apply cryptos list update (price & volume refresh) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.787.59'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.645.82'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.40%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.63'
$ws.Range('E5').Value = '  +0.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('E6').Value = '  -0.64%  '
$ws.Range('E7').Value = '  +0.25%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0628'
$ws.Range('E9').Value = '  +0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.20'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').Value = '1.650.06'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('E13').Value = '  -0.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.528'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '64.68'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('D16').Value = '26.770.43'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('D17').Value = '0.0₃0735'
$ws.Range('E17').Value = '  -1.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '214.34'
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.40'
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('E21').Value = '  +12.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.25'
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.37'
$ws.Range('E23').Value = '  -1.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '146.35'
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.16'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.67'
$ws.Range('E28').Value = '  -0.53%  '
$ws.Range('E29').Value = '  -1.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.37'
$ws.Range('E31').Value = '  -0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.00'
$ws.Range('E32').Value = '  -1.35%  '
$ws.Range('D33').Value = '1.288.62'
$ws.Range('E33').Value = '  +0.72%  '
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.44'
$ws.Range('E35').Value = '  +1.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0177'
$ws.Range('E36').Value = '  -2.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.537'
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.822'
$ws.Range('E38').Value = '  -0.68%  '
$ws.Range('E39').Value = '  +0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.807'
$ws.Range('E40').Value = '  -1.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.23'
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('E42').Value = '  -2.73%  '
$ws.Range('D43').Value = '1.783.16'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.91'
$ws.Range('E44').Value = '  +3.40%  '
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('E47').Value = '  -1.15%  '
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.63'
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0970'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.407'
$ws.Range('E51').Value = '  +0.07%  '
